# Apply updated "想去人数" (interest count) and "最低票价" (min price) figures
# captured at a later crawl timestamp (gh-pages data refresh, commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4587
$ws.Range("F3").Value = 448
$ws.Range("F4").Value = 3699
$ws.Range("G4").Value = 108
$ws.Range("F5").Value = 1082
$ws.Range("F6").Value = 172
$ws.Range("F7").Value = 1347
$ws.Range("F8").Value = 377
$ws.Range("F9").Value = 382
$ws.Range("F10").Value = 2585
$ws.Range("F11").Value = 1293
$ws.Range("F12").Value = 43
$ws.Range("F14").Value = 283
$ws.Range("F16").Value = 565
$ws.Range("F17").Value = 267
$ws.Range("F18").Value = 67
$ws.Range("F19").Value = 10735
$ws.Range("F20").Value = 6180
$ws.Range("F22").Value = 11
$ws.Range("F24").Value = 222
$ws.Range("F27").Value = 852
$ws.Range("F28").Value = 29
$ws.Range("F29").Value = 197
$ws.Range("F30").Value = 870
$ws.Range("F31").Value = 3576
$ws.Range("F32").Value = 49
$ws.Range("F33").Value = 973
$ws.Range("F34").Value = 486
$ws.Range("F35").Value = 138
$ws.Range("F36").Value = 283
$ws.Range("F37").Value = 250
$ws.Range("F39").Value = 4883
$ws.Range("F41").Value = 1162
$ws.Range("F42").Value = 176
$ws.Range("F43").Value = 216
$ws.Range("F44").Value = 126
$ws.Range("F45").Value = 501

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = "不可售"
$ws.Range("F14").Value = 3619

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8877
$ws.Range("F3").Value = 451
$ws.Range("F4").Value = 1689

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8877
$ws.Range("F3").Value = 451
$ws.Range("F4").Value = 1689
$ws.Range("F5").Value = 4587
$ws.Range("F6").Value = 3699
$ws.Range("G6").Value = 108
$ws.Range("F7").Value = 1082
$ws.Range("F8").Value = 172
$ws.Range("F10").Value = 382
$ws.Range("F11").Value = 2585
$ws.Range("F16").Value = 1293
$ws.Range("F18").Value = 43
$ws.Range("F19").Value = 283
$ws.Range("F21").Value = 565
$ws.Range("F22").Value = 267
$ws.Range("F23").Value = 10735
$ws.Range("F24").Value = 3619
$ws.Range("F26").Value = 11
$ws.Range("F28").Value = 222
$ws.Range("F31").Value = 852
$ws.Range("F32").Value = 29
$ws.Range("F33").Value = 870
$ws.Range("F34").Value = 3576
$ws.Range("F35").Value = 49
$ws.Range("F36").Value = 973
$ws.Range("F37").Value = 138
$ws.Range("F38").Value = 283
$ws.Range("F42").Value = 4883
$ws.Range("F44").Value = 1162
$ws.Range("F45").Value = 176
$ws.Range("F46").Value = 126
$ws.Range("F47").Value = 501
